# Update performance dashboard 2025-12-20 10:06
#
# This script rewrites the Pattern1-Pure Data performance rows (for models
# deepseek-v3, gemini-3-pro, gpt-5 and llama-3.1-405b) on both the
# "Summary" sheet and the dedicated "Pattern1-Pure Data" sheet with the
# latest trading-performance figures.
#
# All of the text-like metric columns (currency, percentages, date codes)
# are stored in the workbook as literal text, e.g. "+0.45%" rather than a
# numeric 0.45%. Assigning a string that merely *looks* like a number or a
# percentage to Range.Value causes Excel to auto-convert it into a real
# number, which would change both the stored type and (for percentages)
# introduce a number format that was not present before. To keep such
# values as genuine text - exactly like the source cells - a leading
# apostrophe (the standard Excel "treat as text" quote prefix) is used for
# any replacement value that would otherwise be auto-converted; the
# apostrophe itself is not stored and only the text remains.

function Set-TextValue($ws, $addr, $val) {
    if ($val -match '%' -or $val -match '^[0-9]+$') {
        $ws.Range($addr).Value = "'" + $val
    } else {
        $ws.Range($addr).Value = $val
    }
}

function Update-PerformanceRows($ws) {
    # Row 2: deepseek-v3
    Set-TextValue $ws "D2" "¥1,004,526.00"
    Set-TextValue $ws "E2" "¥+4,526.00"
    Set-TextValue $ws "F2" "+0.45%"
    Set-TextValue $ws "G2" "+76.65%"
    $ws.Range("H2").Value = 19.872
    Set-TextValue $ws "I2" "0.00%"
    Set-TextValue $ws "J2" "100.0%"
    Set-TextValue $ws "K2" "0.2262%"
    Set-TextValue $ws "L2" "0.1807%"

    # Row 3: gemini-3-pro
    Set-TextValue $ws "D3" "¥1,004,601.00"
    Set-TextValue $ws "E3" "¥+4,601.00"
    Set-TextValue $ws "F3" "+0.46%"
    Set-TextValue $ws "G3" "+78.32%"
    $ws.Range("H3").Value = 28.141
    Set-TextValue $ws "K3" "0.2299%"
    Set-TextValue $ws "L3" "0.1297%"
    $ws.Range("M3").Value = 3
    Set-TextValue $ws "O3" "20251219"

    # Row 4: gpt-5
    Set-TextValue $ws "D4" "¥1,003,469.00"
    Set-TextValue $ws "E4" "¥+3,469.00"
    Set-TextValue $ws "F4" "+0.35%"
    Set-TextValue $ws "G4" "+54.70%"
    $ws.Range("H4").Value = 21.573
    Set-TextValue $ws "I4" "0.00%"
    Set-TextValue $ws "J4" "100.0%"
    Set-TextValue $ws "K4" "0.1734%"
    Set-TextValue $ws "L4" "0.1276%"

    # Row 5: llama-3.1-405b
    Set-TextValue $ws "D5" "¥1,001,074.00"
    Set-TextValue $ws "E5" "¥+1,074.00"
    Set-TextValue $ws "F5" "+0.11%"
    Set-TextValue $ws "G5" "+14.48%"
    $ws.Range("H5").Value = 9.238
    Set-TextValue $ws "I5" "0.04%"
    Set-TextValue $ws "J5" "50.0%"
    Set-TextValue $ws "K5" "0.0537%"
    Set-TextValue $ws "L5" "0.0923%"
}

$wb = $excel.ActiveWorkbook

Update-PerformanceRows $wb.Worksheets.Item("Summary")
Update-PerformanceRows $wb.Worksheets.Item("Pattern1-Pure Data")
